# ---------------------------------------------------------------------------
# Add "Cart & Checkout" worksheet (between "Test Cases of Add to Cart" and
# "Bug Report for Login"), populate it with 5 new cart/checkout test cases,
# and tweak the selections/row-height left behind on the two sheets the
# author had open while editing.
# ---------------------------------------------------------------------------

function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$xlPasteFormats = -4122
$xlLeft         = -4131
$xlCenter       = -4108
$xlGeneral      = 1
$xlBottom       = -4107

$wb = $excel.ActiveWorkbook

$loginSheet   = $wb.Worksheets.Item("Test Cases for Login ")
$addCartSheet = $wb.Worksheets.Item("Test Cases of Add to Cart")
$bugSheet     = $wb.Worksheets.Item("Bug Report for Login")

# ---------------------------------------------------------------------------
# 1. Tidy up the two existing sheets the way the author left them.
# ---------------------------------------------------------------------------

# "Test Cases for Login " - just a new selected cell.
$loginSheet.Activate() | Out-Null
$loginSheet.Range("G2").Select() | Out-Null

# "Test Cases of Add to Cart" - header row shrinks back to a normal height,
# and the selection moves to F2 (it also stops being the active/tabbed sheet
# once we activate the new sheet further down).
$addCartSheet.Activate() | Out-Null
$addCartSheet.Rows.Item(1).RowHeight = 18.75
$addCartSheet.Range("F2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert the new sheet right after "Test Cases of Add to Cart".
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $addCartSheet)
$ws.Name = "Cart & Checkout"
$ws.Tab.Color = RGB 255 192 0

# ---------------------------------------------------------------------------
# 3. Header row - copy the look of the other sheets' header row.
# ---------------------------------------------------------------------------

$ws.Range("A1").Value = "Test Case ID"
$ws.Range("B1").Value = "Test Scenario"
$ws.Range("C1").Value = "Precondition"
$ws.Range("D1").Value = "Steps"
$ws.Range("E1").Value = "Expected Result"
$ws.Range("F1").Value = "Priority"
$ws.Range("G1").Value = "Status"

$addCartSheet.Range("A1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("B1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("C1:G1").Copy() | Out-Null
$ws.Range("C1:G1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Rows.Item(1).RowHeight = 18.75

# ---------------------------------------------------------------------------
# 4. Data rows (TC012 .. TC016).
# ---------------------------------------------------------------------------

# Row 2 - TC012
$ws.Range("A2").Value = "TC012"
$ws.Range("B2").Value = "Update item quantity in cart "
$ws.Range("C2").Value = "At least one item in cart"
$ws.Range("D2").Value = "1-Go to Cart `n2-Click on quantity selector for an item`n3-Update the quantity (e.g, from 1 to 3)"
$ws.Range("E2").Value = "Item quantity updates correctly `nTotal price reflect the update quantity "
$ws.Range("F2").Value = "high"
$ws.Range("G2").Value = "Pass"

$addCartSheet.Range("A2:D2").Copy() | Out-Null
$ws.Range("A2:D2").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("E2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("F2:G2").Copy() | Out-Null
$ws.Range("F2:G2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Rows.Item(2).RowHeight = 65.25

# Row 3 - TC013
$ws.Range("A3").Value = "TC013"
$ws.Range("B3").Value = "Checkout/Place Order"
$ws.Range("C3").Value = "At least one item in cart"
$ws.Range("D3").Value = "1-Go to Cart`n2-Click Checkout `n3-Enter shipping information `n4-Enter payment details`n5-Click Place Order"
$ws.Range("E3").Value = "Order is successfully placed Confirmation message is displayed .`nTotal price updated accordingly."
$ws.Range("F3").Value = "high"
$ws.Range("G3").Value = "Pass"

$addCartSheet.Range("A2:D2").Copy() | Out-Null
$ws.Range("A3:D3").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("E2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("F2:G2").Copy() | Out-Null
$ws.Range("F3:G3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Rows.Item(3).RowHeight = 75

# Row 4 - TC014
$ws.Range("A4").Value = "TC014"
$ws.Range("B4").Value = "Apply discount code"
$ws.Range("C4").Value = "At least one item in cart"
$ws.Range("D4").Value = "1-Go to Cart`n2-Enter Valid discount code `n3-Click Apply"
$ws.Range("E4").Value = "Discount applied correctly.`nTotal price updated accordingly."
$ws.Range("F4").Value = "Medium"
$ws.Range("G4").Value = "Pass"

$addCartSheet.Range("A2:D2").Copy() | Out-Null
$ws.Range("A4:D4").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("E2").Copy() | Out-Null
$ws.Range("E4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E4").HorizontalAlignment = $xlGeneral
$loginSheet.Range("F4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("G2").Copy() | Out-Null
$ws.Range("G4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Rows.Item(4).RowHeight = 45

# Row 5 - TC015
$ws.Range("A5").Value = "TC015"
$ws.Range("B5").Value = "Empty Cart"
$ws.Range("C5").Value = "At least one item in cart"
$ws.Range("D5").Value = "1-Click to Cart`n2-Click `"Empty Cart`"button `n3-Confirm action"
$ws.Range("E5").Value = "Cart is emptied `nCart badge show 0"
$ws.Range("F5").Value = "Medium"
$ws.Range("G5").Value = "Pass"

$addCartSheet.Range("A2:D2").Copy() | Out-Null
$ws.Range("A5:D5").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("E2").Copy() | Out-Null
$ws.Range("E5").PasteSpecial($xlPasteFormats) | Out-Null
$loginSheet.Range("F4").Copy() | Out-Null
$ws.Range("F5").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("G2").Copy() | Out-Null
$ws.Range("G5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Rows.Item(5).RowHeight = 45

# Row 6 - TC016
$ws.Range("A6").Value = "TC016"
$ws.Range("B6").Value = "View Cart Details"
$ws.Range("C6").Value = "At least one item in cart"
$ws.Range("D6").Value = "1-Go to Cart"
$ws.Range("E6").Value = "All items,quantity and total price are displayed correctly"
$ws.Range("F6").Value = "Low"
$ws.Range("G6").Value = "Pass"

$addCartSheet.Range("A2:D2").Copy() | Out-Null
$ws.Range("A6:D6").PasteSpecial($xlPasteFormats) | Out-Null
$addCartSheet.Range("E2").Copy() | Out-Null
$ws.Range("E6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E6").HorizontalAlignment = $xlGeneral
$ws.Range("E6").VerticalAlignment = $xlBottom
$addCartSheet.Range("F2").Copy() | Out-Null
$ws.Range("F6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F6").Font.Size = 14
$addCartSheet.Range("G2").Copy() | Out-Null
$ws.Range("G6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Rows.Item(6).RowHeight = 30

# ---------------------------------------------------------------------------
# 5. Column widths, sheet layout, and final selection/active sheet.
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 14.307292
$ws.Columns.Item(2).ColumnWidth = 26.022135
$ws.Columns.Item(3).ColumnWidth = 21.592448
$ws.Columns.Item(4).ColumnWidth = 35.877604
$ws.Columns.Item(5).ColumnWidth = 35.166667
$ws.Columns.Item(6).ColumnWidth = 12.307292
$ws.Columns.Item(7).ColumnWidth = 7.592448

$ws.Range("A1:G1").Borders.LineStyle = 1

$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null
